$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Put the e-mail address text in A2, then turn it into a real hyperlink.
# Assigning the text first (rather than relying on Hyperlinks.Add's
# TextToDisplay parameter) keeps the cell's display text exactly equal
# to the address, matching what Excel does when you type an e-mail
# address into a cell and it auto-converts to a hyperlink.
$ws.Range("A2").Value = "hermanliran@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:hermanliran@gmail.com")
